$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Widen column B (used for long descriptive text) to match the author's
# resize. Excel's COM ColumnWidth is expressed in "characters" of the
# workbook's default font and gets re-quantised to pixels on save, so the
# raw OOXML <col width> that ends up on disk isn't a 1:1 echo of the value
# assigned here. 71.17 is the character-width input that this runtime's
# pixel/font metrics round-trip to an OOXML width of exactly 72, which is
# the value seen in the target file.
$ws.Columns.Item(2).ColumnWidth = 71.17
